$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "28.489.59"
$ws.Range("E2").Value = "  -3.70%  "
Set-TextValue "D3" "1.947.84"
$ws.Range("E3").Value = "  -2.79%  "
Set-TextValue "D4" "1.014"
$ws.Range("E4").Value = "  +0.14%  "
Set-TextValue "D5" "321.54"
$ws.Range("E5").Value = "  -2.55%  "
$ws.Range("E6").Value = "  +0.23%  "
Set-TextValue "D7" "0.4764"
$ws.Range("E7").Value = "  -5.03%  "
$ws.Range("E8").Value = "  -4.88%  "
Set-TextValue "D9" "53.61"
$ws.Range("E9").Value = "  -0.77%  "
Set-TextValue "D10" "0.08544"
$ws.Range("E10").Value = "  -5.15%  "
Set-TextValue "D11" "1.060"
$ws.Range("E11").Value = "  -5.27%  "
Set-TextValue "D12" "22.03"
$ws.Range("E12").Value = "  -5.44%  "
Set-TextValue "D13" "1.957.89"
$ws.Range("E13").Value = "  -1.88%  "
Set-TextValue "D14" "7.596"
$ws.Range("E14").Value = "  -5.66%  "
Set-TextValue "D15" "6.175"
$ws.Range("E15").Value = "  -4.88%  "
Set-TextValue "D16" "1.015"
$ws.Range("E16").Value = "  +0.11%  "
Set-TextValue "D17" "0.00001078"
$ws.Range("E17").Value = "  -3.19%  "
Set-TextValue "D18" "88.86"
$ws.Range("E18").Value = "  -5.90%  "
Set-TextValue "D19" "0.06635"
$ws.Range("E19").Value = "  -0.80%  "
Set-TextValue "D20" "18.68"
$ws.Range("E20").Value = "  -5.17%  "
Set-TextValue "D21" "1.014"
$ws.Range("E21").Value = "  +0.34%  "
Set-TextValue "D22" "5.806"
$ws.Range("E22").Value = "  -2.67%  "
Set-TextValue "D23" "28.489.63"
$ws.Range("E23").Value = "  -3.80%  "
Set-TextValue "D24" "11.49"
$ws.Range("E24").Value = "  -4.35%  "
$ws.Range("E25").Value = "  -0.38%  "
Set-TextValue "D26" "2.183.48"
$ws.Range("E26").Value = "  -2.28%  "
Set-TextValue "D27" "153.90"
$ws.Range("E27").Value = "  -3.17%  "
Set-TextValue "D28" "20.16"
$ws.Range("E28").Value = "  -2.75%  "
Set-TextValue "D29" "5.932"
$ws.Range("E29").Value = "  -7.34%  "
Set-TextValue "D30" "2.157"
$ws.Range("E30").Value = "  -6.20%  "
Set-TextValue "D31" "123.52"
$ws.Range("E31").Value = "  -3.71%  "
Set-TextValue "D32" "0.9928"
$ws.Range("E32").Value = "  -6.18%  "
Set-TextValue "D33" "0.09541"
$ws.Range("E33").Value = "  -4.28%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D34" "1.446"
$ws.Range("E34").Value = "  -7.85%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D35" "3.671"
$ws.Range("E35").Value = "  -3.40%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D36" "5.588"
$ws.Range("E36").Value = "  -4.37%  "
Set-TextValue "D37" "0.02336"
$ws.Range("E37").Value = "  -5.65%  "
Set-TextValue "D38" "0.06217"
$ws.Range("E38").Value = "  -2.77%  "
Set-TextValue "D39" "1.259"
$ws.Range("E39").Value = "  -3.86%  "
Set-TextValue "D40" "8.729"
$ws.Range("E40").Value = "  -6.19%  "
$ws.Range("E41").Value = "  -5.17%  "
Set-TextValue "D42" "11.07"
$ws.Range("E42").Value = "  -5.47%  "
Set-TextValue "D43" "1.013"
$ws.Range("E43").Value = "  +0.22%  "
Set-TextValue "D44" "0.1925"
$ws.Range("E44").Value = "  -6.42%  "
Set-TextValue "D45" "1.330"
$ws.Range("E45").Value = "  +1.96%  "
Set-TextValue "D46" "0.5951"
$ws.Range("E46").Value = "  -6.43%  "
Set-TextValue "D47" "12.98"
$ws.Range("E47").Value = "  -3.28%  "
Set-TextValue "D48" "2.061"
$ws.Range("E48").Value = "  -6.14%  "
Set-TextValue "D49" "3.397"
$ws.Range("E49").Value = "  -3.35%  "
Set-TextValue "D50" "0.00000000332"
$ws.Range("E50").Value = "  -0.62%  "
Set-TextValue "D51" "0.06802"
$ws.Range("E51").Value = "  -2.70%  "
